# Apply the receipts workbook edit:
#  - rename sheet "Sheet1" -> "2025-04-23"
#  - fill in the (previously blank) header row 1
#  - correct/update a few cell values in the existing rows 2 and 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab.
$ws.Name = "2025-04-23"

# Row 1 was reserved (dimension already spanned A1:H3) but empty; populate it
# with the column headers.
$ws.Range("A1").Value = "Sponsor Name"
$ws.Range("B1").Value = "Guest Name"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Initials"
$ws.Range("E1").Value = "Receipt Number"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Item"
$ws.Range("H1").Value = "Price"

# Row 2 (Matthew Wolz / JJ  / 4/23/2025 / MW / 1 / N/A / Daily Guest Pass / 3)
# The date is entered with a leading apostrophe so it is kept as literal text
# ("2025-04-23") instead of being parsed into a date serial number.
$ws.Range("B2").Value = "Tiffany"
$ws.Range("C2").Value = "'2025-04-23"
$ws.Range("E2").Value = 6

# Row 3 (Matthew wolz / JJ / 4/23/2025 / MW / 1 / N/A / Daily Guest Pass / 3)
$ws.Range("A3").Value = "Matthew Wolz"
$ws.Range("B3").Value = "JJ James JJ"
$ws.Range("C3").Value = "'2025-04-23"
$ws.Range("E3").Value = 7
